$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: lower-case the labels and rename the metric columns ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data rows: columns D/E/F are reshuffled (D<-old E, E<-old F) and F gets new climate-change values ---
$ws.Range("D2").Value = 1.445327933333333
$ws.Range("E2").Value = 22.026247
$ws.Range("F2").Value = 0.000040299577

$ws.Range("D3").Value = 1.818808266666667
$ws.Range("E3").Value = 23.2696
$ws.Range("F3").Value = 0.0000507132

$ws.Range("D4").Value = 2.8906558
$ws.Range("E4").Value = 44.052493
$ws.Range("F4").Value = 0.000080599154

$ws.Range("D5").Value = 1.1567294
$ws.Range("E5").Value = 17.792089
$ws.Range("F5").Value = 0.000032252684

$ws.Range("D6").Value = 2.602057266666667
$ws.Range("E6").Value = 39.818336
$ws.Range("F6").Value = 0.00007255226100000001

$ws.Range("D7").Value = 2.8906558
$ws.Range("E7").Value = 44.052493
$ws.Range("F7").Value = 0.000080599154

$ws.Range("D8").Value = 1.016215933333333
$ws.Range("E8").Value = 15.379507
$ws.Range("F8").Value = 0.000028334796

$ws.Range("D9").Value = 1.127785
$ws.Range("E9").Value = 17.10766
$ws.Range("F9").Value = 0.000031445639

$ws.Range("D10").Value = 1.445327933333333
$ws.Range("E10").Value = 22.026247
$ws.Range("F10").Value = 0.000040299577

$ws.Range("D11").Value = 1.273683133333333
$ws.Range("E11").Value = 19.367551
$ws.Range("F11").Value = 0.000035513664

$ws.Range("D12").Value = 1.230771866666667
$ws.Range("E12").Value = 18.702877
$ws.Range("F12").Value = 0.000034317186

$ws.Range("D13").Value = 1.102038333333333
$ws.Range("E13").Value = 16.708855
$ws.Range("F13").Value = 0.000030727752

$ws.Range("D14").Value = 1.222189666666667
$ws.Range("E14").Value = 18.569942
$ws.Range("F14").Value = 0.000034077891

# --- Add header-row cell comments describing each column's data type ---
$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null
